# Applies the attendance_reports sync edits described in the commit:
# "Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-08 07:19:51"
#
# The edits are:
#  - Re-order several "Recorded By" email lists (same members, new order)
#  - Update a handful of numeric counters (Missing/Pending sessions, etc.)
#  - Update the attendance fraction for one session (71/251 -> 72/251)
#  - Flip row 25 (PATHOLOGY LAB/MUSEUM session on 08/12/2025) from the
#    "Pending" (yellow) state to "Not Recorded" (pink) state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-order "Recorded By" (column G) email lists -------------------------

$ws.Range("G2").Value = "gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System, servinaz@med.asu.edu.eg"

$ws.Range("G3").Value = "majorelle.magdy@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, System"

$ws.Range("G4").Value = "asmaa.reda@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, gehanadel@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg"

$ws.Range("G5").Value = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"

$ws.Range("G6").Value = "majorelle.magdy@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, alshimaa.atef@med.asu.edu.egm"

$ws.Range("G7").Value = "AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, Amera.a.saad@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Kerelos.zareef@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg"

$ws.Range("G12").Value = "dina.adel@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg"

$ws.Range("G27").Value = "hana.amr@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"

$ws.Range("G30").Value = "yassmen.ahmed@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg"

# --- Numeric counter updates ------------------------------------------------

$ws.Range("L7").Value = 3
$ws.Range("L8").Value = 4

$ws.Range("P15").Value = 3
$ws.Range("Q15").Value = 4

# --- Attendance fraction update --------------------------------------------

$ws.Range("H23").Value = "72/251"

# --- Row 25 status flip: Pending (yellow) -> Not Recorded (pink) -----------
# Row 11 already uses the "Not Recorded" pink formatting; copy its formats
# across to row 25 so the fill/font match exactly, then update the text.

$ws.Range("A11:I11").Copy()
$ws.Range("A25:I25").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I25").Value = "Not Recorded"
